$wb = $excel.ActiveWorkbook

$oldDate = "January 30 2026 16.19.47 EST"
$newDate = "February 02 2026 12.49.33 EST"

$aboutWs = $wb.Worksheets.Item("About")
$boundariesWs = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the "About" sheet cells that contain the build version/date.
$cellsToUpdate = @("A2", "A6")
foreach ($addr in $cellsToUpdate) {
    $cell = $aboutWs.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null -and $text.Contains($oldDate)) {
        $cell.Value = $text.Replace($oldDate, $newDate)
    }
}

# Update the build_version column (S) on the "Boundaries and methane sources" sheet.
for ($row = 2; $row -le 10; $row++) {
    $cell = $boundariesWs.Range("S" + $row)
    $text = $cell.Value2
    if ($text -ne $null -and $text.Contains($oldDate)) {
        $cell.Value = $text.Replace($oldDate, $newDate)
    }
}
